# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect newly-scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 731
$ws1.Range("F4").Value  = 1440
$ws1.Range("F8").Value  = 6116
$ws1.Range("F12").Value = 4965
$ws1.Range("F13").Value = 24
$ws1.Range("F15").Value = 1166
$ws1.Range("F17").Value = 353
$ws1.Range("F18").Value = 56
$ws1.Range("F22").Value = 3462

# Sheet "全部类型" (sheet4) - same events, one row lower
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 731
$ws4.Range("F5").Value  = 1440
$ws4.Range("F9").Value  = 6116
$ws4.Range("F13").Value = 4965
$ws4.Range("F14").Value = 24
$ws4.Range("F16").Value = 1166
$ws4.Range("F18").Value = 353
$ws4.Range("F19").Value = 56
$ws4.Range("F23").Value = 3462
